$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 1176
$ws1.Range("F5").Value = 19
$ws1.Range("F7").Value = 288
$ws1.Range("F8").Value = 389
$ws1.Range("F9").Value = 1031
$ws1.Range("F14").Value = 13190
$ws1.Range("F18").Value = 5424
$ws1.Range("F19").Value = 5560
$ws1.Range("F20").Value = 25

# Sheet "全部类型" (sheet4) updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1176
$ws4.Range("F12").Value = 19
$ws4.Range("F23").Value = 288
$ws4.Range("F24").Value = 389
$ws4.Range("F31").Value = 1031
$ws4.Range("F36").Value = 13190
$ws4.Range("F41").Value = 5424
$ws4.Range("F42").Value = 5560
$ws4.Range("F43").Value = 25
